$wb = $excel.ActiveWorkbook
$tc1 = $wb.Worksheets.Item("tc1")

$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $tc1)
$ws.Name = "ValidLogin"

$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"
$ws.Range("B2").Value = "pointofsale"
$ws.Range("A2").Value = "ADMIN"

$ws.Activate()
$excel.ActiveWindow.Zoom = 205
$ws.Range("A2").Select() | Out-Null
